$d = $word.ActiveDocument

# --- 1. Remove the second run's text (" lớp 12522T.1") entirely -----------
# This also removes the run element itself once its text is empty.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(" lớp 12522T.1", $false, $false, $false, $false, $false, `
              $true, 1, $true, "", 2)

# --- 2. Fix the casing of the remaining word: "Test" -> "test" ------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("Test", $true, $false, $false, $false, $false, `
               $true, 1, $true, "test", 2)

# --- 3. Restore the paragraph/run language to English (US) ----------------
# Re-apply en-US across the whole (now single-run) paragraph so every run's
# w:lang matches the restored original.
$d.Content.LanguageID = "en-US"
